# Update countries & provincias Spain
# Applies a refreshed-data update to the "Pais" sheet:
#  - updates the "last updated" timestamp label
#  - updates case-count figures for several countries (refreshed source data)
#  - Bahamas / Monaco swap rank (row 181 becomes Bahamas, row 182 becomes Monaco)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Timestamp label in A1
$ws.Range("A1").Value = "Datos actualizados a 11 de Julio de 2020 a las 04:55"

# Estados Unidos row (row 5) - Casos activos / Recuperados refreshed
$ws.Range("D5").Value = 1213512
$ws.Range("E5").Value = 520302

# Row 43 refreshed figures
$ws.Range("B43").Value = 45565
$ws.Range("C43").Value = 1452
$ws.Range("D43").Value = 13918
$ws.Range("E43").Value = 29945
$ws.Range("G43").Value = 64
$ws.Range("H43").Value = 1702

# Rows 181/182: Bahamas and Monaco swap rank with refreshed figures
$ws.Range("A181").Value = "Bahamas"
$ws.Range("B181").Value = 108
$ws.Range("C181").Value = 1
$ws.Range("D181").Value = 89
$ws.Range("E181").Value = 8
$ws.Range("F181").Value = 0
$ws.Range("G181").Value = 0
$ws.Range("H181").Value = 11

$ws.Range("A182").Value = "Monaco"
$ws.Range("B182").Value = 108
$ws.Range("C182").Value = 0
$ws.Range("D182").Value = 96
$ws.Range("E182").Value = 8
$ws.Range("F182").Value = 0
$ws.Range("G182").Value = 0
$ws.Range("H182").Value = 4

# Row 196 refreshed figures
$ws.Range("B196").Value = 37
$ws.Range("C196").Value = 4
$ws.Range("E196").Value = 15
